# Update the "K" column (column G) with recomputed strike values.
# Commit message: regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 0
    6  = 3
    8  = 2
    10 = 3
    11 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
